$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "ITGC-OP-01"
$ws.Range("B2").Value = "테스트 접근권한 관리"
$ws.Range("C2").Value = "시스템 접근 권한을 적절히 부여하고 관리한다."
$ws.Range("E2").Value = "상시"
$ws.Range("H2").Value = "접근권한 목록"
$ws.Range("I2").Value = "권한 부여 현황 확인"

# Update row 3
$ws.Range("A3").Value = "ITGC-OP-02"
$ws.Range("B3").Value = "테스트 변경관리"
$ws.Range("C3").Value = "시스템 변경 시 승인 절차를 따른다."
$ws.Range("D3").Value = "Y"
$ws.Range("E3").Value = "수시"
$ws.Range("F3").Value = "탐지"
$ws.Range("H3").Value = "변경요청서"
$ws.Range("I3").Value = "변경 승인 이력 확인"

# Add new row 4
$ws.Range("A4").Value = "ITGC-OP-03"
$ws.Range("B4").Value = "테스트 운영 보안"
$ws.Range("C4").Value = "운영 환경의 보안을 유지한다."
$ws.Range("D4").Value = "N"
$ws.Range("E4").Value = "월별"
$ws.Range("F4").Value = "예방"
$ws.Range("G4").Value = "수동"
$ws.Range("H4").Value = "보안점검표"
$ws.Range("I4").Value = "월별 점검 결과 확인"
